$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last-updated timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 16:51"

# Country label swaps (ranking shuffled between refreshes)
$ws.Range("A28").Value = "Irak"
$ws.Range("A29").Value = "Indonesia"
$ws.Range("A44").Value = "Republica Dominicana"
$ws.Range("A45").Value = "Panama"

# Updated statistics
# Row 4
$ws.Range("B4").Value = 3358120
$ws.Range("C4").Value = 2474
$ws.Range("D4").Value = 1490724
$ws.Range("E4").Value = 1729967
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 137429
# Row 5
$ws.Range("B5").Value = 1842127
$ws.Range("C5").Value = 1315
$ws.Range("E5").Value = 557100
$ws.Range("G5").Value = 23
$ws.Range("H5").Value = 71515
# Row 6
$ws.Range("B6").Value = 867020
$ws.Range("C6").Value = 16662
$ws.Range("E6").Value = 303031
# Row 25
$ws.Range("D25").Value = 42694
$ws.Range("E25").Value = 52997
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 1818
# Row 28
$ws.Range("B28").Value = 77506
$ws.Range("C28").Value = 2312
$ws.Range("D28").Value = 44724
$ws.Range("E28").Value = 29632
$ws.Range("G28").Value = 95
$ws.Range("H28").Value = 3150
# Row 29
$ws.Range("B29").Value = 75699
$ws.Range("C29").Value = 1681
$ws.Range("D29").Value = 35638
$ws.Range("E29").Value = 36455
$ws.Range("G29").Value = 71
$ws.Range("H29").Value = 3606
# Row 44
$ws.Range("B44").Value = 44532
$ws.Range("C44").Value = 1418
$ws.Range("D44").Value = 21459
$ws.Range("E44").Value = 22176
$ws.Range("G44").Value = 17
$ws.Range("H44").Value = 897
# Row 45
$ws.Range("B45").Value = 44332
$ws.Range("D45").Value = 22170
$ws.Range("E45").Value = 21269
$ws.Range("H45").Value = 893
# Row 58
$ws.Range("B58").Value = 24041
$ws.Range("C58").Value = 520
$ws.Range("D58").Value = 15093
$ws.Range("E58").Value = 8642
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 306
# Row 60
$ws.Range("B60").Value = 19382
$ws.Range("C60").Value = 174
$ws.Range("E60").Value = 6073
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 642
# Row 95
$ws.Range("B95").Value = 4972
$ws.Range("C95").Value = 4
$ws.Range("D95").Value = 4712
$ws.Range("E95").Value = 204
# Row 112
$ws.Range("B112").Value = 2411
$ws.Range("C112").Value = 5
$ws.Range("D112").Value = 1694
$ws.Range("E112").Value = 596
# Row 128
$ws.Range("B128").Value = 1433
$ws.Range("C128").Value = 44
$ws.Range("D128").Value = 341
$ws.Range("E128").Value = 1053
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 39
# Row 164
$ws.Range("B164").Value = 331
$ws.Range("C164").Value = 1
$ws.Range("D164").Value = 261
